# Insert a new column "statut_name" at column C, shifting the existing
# NCTId..intervention_type columns (previously C..L) one position to the
# right (now D..M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:L one position to the right to make room for the new column.
$ws.Range("C:C").EntireColumn.Insert()

# Header for the newly inserted column, matching the bold/centered/bordered
# style used by the rest of row 1 (copy formatting from the neighbouring
# header cell that was just shifted into D1).
$ws.Range("C1").Value = "statut_name"
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new column for each data row with the constant status text.
$statutName = "pas de résultat ni de publication"
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = $statutName
}
